# Add a new worksheet "V7" by duplicating "V6" (same layout/tables/styles),
# then update it with the new model run's data (per commit message:
# "serialized version 7 of model with more data. 79% phishing, 89% accuracy overall").

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("V6")
$src.Copy($null, $src)
$ws = $wb.Worksheets.Item("V6 (2)")
$ws.Name = "V7"

# ---- Non phish mailbox block (rows 2-8) ----
$ws.Range("A5").Value = "Joy_Mailbox"
$ws.Range("B5").Value = 79

$ws.Range("A6").Value = "Yannis_Mailbox"
$ws.Range("B6").Value = 182

$ws.Range("A7").Value = "ModernHam1"
$ws.Range("B7").Value = 88

$ws.Range("B8").Formula = "=SUM(B2:B7)"

# ---- Modern Day Ham test table (D3:G8 / rows 4-8) ----
$ws.Range("D5").Value = 2
$ws.Range("D6").Value = 3

$ws.Range("E4").Value = 222
$ws.Range("F4").Value = 222
$ws.Range("E5").Value = 222
$ws.Range("F5").Value = 222
$ws.Range("E6").Value = 222
$ws.Range("F6").Value = 222
$ws.Range("E7").Value = 222
$ws.Range("F7").Value = 222
$ws.Range("E8").Value = 222
$ws.Range("F8").Value = 222

# ---- Phish mailbox block (rows 11-13) ----
$ws.Range("A13").Value = "ModernPhish2"
$ws.Range("B13").Value = 56

# ---- Modern Day Phish test table (D13:G18 / rows 14-18) ----
$ws.Range("E14").Value = 46
$ws.Range("F14").Value = 56
$ws.Range("E15").Value = 44
$ws.Range("F15").Value = 56
$ws.Range("E16").Value = 42
$ws.Range("F16").Value = 56
$ws.Range("E17").Value = 44
$ws.Range("F17").Value = 56
$ws.Range("E18").Value = 47
$ws.Range("F18").Value = 56

# ---- Model metadata ----
$ws.Range("B21").Value = "dsv7.csv"

# ---- Conclusion ----
$ws.Range("A25").Value = "Added about 40 mote items to train and increased phishing sample size. Model sits at about 79.6% right now, conclusion: requires more phishing emails"
$ws.Rows.Item(25).RowHeight = 75

# ---- Rebuild the 4 ListObjects/tables that a straight sheet-copy doesn't carry over ----
$t1 = $ws.ListObjects.Add(1, $ws.Range("D3:G8"), $null, 1)
$t1.Name = "Table411281822"
$t1.TableStyle = "TableStyleLight15"

$t2 = $ws.ListObjects.Add(1, $ws.Range("D13:G18"), $null, 1)
$t2.Name = "Table46123151923"
$t2.TableStyle = "TableStyleLight15"

$t3 = $ws.ListObjects.Add(1, $ws.Range("I3:L8"), $null, 1)
$t3.Name = "Table49134162024"
$t3.TableStyle = "TableStyleLight15"

$t4 = $ws.ListObjects.Add(1, $ws.Range("I13:L18"), $null, 1)
$t4.Name = "Table4610147172125"
$t4.TableStyle = "TableStyleLight15"

# ---- Selection / active sheet ----
$ws.Range("H16").Select()
$ws.Activate()
